$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5 from 2023-11-03 (45233) to 2023-11-13 (45243)
$ws.Range("C2").Value = 45243
$ws.Range("C3").Value = 45243
$ws.Range("C4").Value = 45243
$ws.Range("C5").Value = 45243
